$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "VALOR MORA" total (doubled once the new period is added)
$ws.Range("E11").Value = 341640

# Update "Cant. Periodos" count (now 2 periods: 2507 and 2508)
$ws.Range("F13").Value = 2

# Insert 3 new data rows before the closing/total row (old row 18),
# pushing it down to row 21 along with everything after it.
$ws.Rows("18:20").Insert()

# Copy the formatting used by the existing data rows (16/17) onto
# the three freshly inserted rows so borders/fonts/fills match.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 18: JONATAN DAVID GRACIA MARRUGO, period 2507 (moved up from the old row 18)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047409050"
$ws.Range("D18").Value = "JONATAN DAVID GRACIA MARRUGO"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: PAOLA KARINA MEÑACA RUIZ, new period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "45549762"
$ws.Range("D19").Value = "PAOLA KARINA MEÑACA RUIZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20: ARGEMIRO VALENCIA MERCADO, new period 2508
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1002258238"
$ws.Range("D20").Value = "ARGEMIRO VALENCIA MERCADO"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# Row 21 (originally row 18, shifted down by the insert): JONATAN DAVID
# GRACIA MARRUGO now reported for the new period 2508.
$ws.Range("E21").Value = "2508"
